# Apply updates described in the commit "Update PPR and brucellosis AHLE
# parameter files" to the PPR_AHLE SMALLRUMINANTS scenario-parameter sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 15 ("part" / Fertility) — refreshed rgamma fits for K:N ---
$ws.Range("K15").Value = "rgamma(10000,9.8731995367749,11.2128609319927)"
$ws.Range("L15").Value = "rgamma(10000,45.545556222724,56.8170215196901)"
$ws.Range("M15").Value = "rgamma(10000,19.021394202597,23.6288060522362)"
$ws.Range("N15").Value = "rgamma(10000,61.0324282301941,72.345654204045)"

# --- New row 29 : single value added in column G ---
$ws.Range("G29").Value = "'8.1666666666666665E-2"

# --- Row 30 header row gains numeric values in B and G ---
$ws.Range("B30").Value = "'1.4999999999999999E-2"
$ws.Range("G30").Value = "'3.3333333333333333E-2"

# --- Row 31 "AlphaN" — refreshed rgamma fits for K:N ---
$ws.Range("K31").Value = "rgamma(10000,3.61170067257074,3403.43585425309)"
$ws.Range("L31").Value = "rgamma(10000,3.60960309269834,3234.05400729713)"
$ws.Range("M31").Value = "rgamma(10000,2.89973650608414,1382.80543309279)"
$ws.Range("N31").Value = "rgamma(10000,2.96996945220457,1352.89487200589)"

# --- Row 32 "AlphaJ" — refreshed rgamma fits for K:N ---
$ws.Range("K32").Value = "rgamma(10000,4.21154853218239,3539.49622589056)"
$ws.Range("L32").Value = "rgamma(10000,4.19063547534716,3306.6993145909)"
$ws.Range("M32").Value = "rgamma(10000,3.15274454073964,1465.75890925845)"
$ws.Range("N32").Value = "rgamma(10000,3.16041418729247,1413.90407390784)"

# --- Row 33 "AlphaF" — refreshed rgamma fits for K:N ---
$ws.Range("K33").Value = "rgamma(10000,4.44815194346767,4002.40107198309)"
$ws.Range("L33").Value = "rgamma(10000,4.60386165488748,3894.74739558975)"
$ws.Range("M33").Value = "rgamma(10000,3.4930659319298,1602.88553751651)"
$ws.Range("N33").Value = "rgamma(10000,3.54212155682475,1564.95658431496)"

# --- Row 34 "AlphaM" — refreshed rgamma fits for K:N ---
$ws.Range("K34").Value = "rgamma(10000,4.46975086236396,4042.30975498851)"
$ws.Range("L34").Value = "rgamma(10000,4.60818773709151,3904.21069711548)"
$ws.Range("M34").Value = "rgamma(10000,3.49084487076966,1613.46183712223)"
$ws.Range("N34").Value = "rgamma(10000,3.46785111439653,1526.74941595658)"

# --- Row 115 "Health_exp" — refreshed rgamma fits for K:N ---
$ws.Range("K115").Value = "rgamma(10000,7.93941338084862,272.859442324841)"
$ws.Range("L115").Value = "rgamma(10000,7.76882162737506,250.445061203992)"
$ws.Range("M115").Value = "rgamma(10000,9.19838639987883,256.467213555895)"
$ws.Range("N115").Value = "rgamma(10000,9.22030903876569,247.35385289858)"
